$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns: G (lockdown_freq) and I (test_freq, after new reward col)
$ws.Range("G1").EntireColumn.Insert() | Out-Null
$ws.Range("I1").EntireColumn.Insert() | Out-Null

# Set header row
$ws.Range("B1").Value = "deaths"
$ws.Range("C1").Value = "delta_schooling"
$ws.Range("D1").Value = "economics_value"
$ws.Range("E1").Value = "heuristic"
$ws.Range("F1").Value = "icus"
$ws.Range("G1").Value = "lockdown_freq"
$ws.Range("H1").Value = "reward"
$ws.Range("I1").Value = "test_freq"
$ws.Range("J1").Value = "testing"
$ws.Range("K1").Value = "tests"
$ws.Range("L1").Value = "xi"

# Set data rows
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1549.820285967816
$ws.Range("C2").Value = 0.5
$ws.Range("D2").Value = 43844049464.29655
$ws.Range("E2").Value = "linearization_heuristic"
$ws.Range("F2").Value = 2000
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 40250755195.01114
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = "linearization_heuristic"
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2231941.8

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1577.388424166159
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 43845881124.38062
$ws.Range("E3").Value = "linearization_heuristic"
$ws.Range("F3").Value = 2000
$ws.Range("G3").Value = 14
$ws.Range("H3").Value = 40185246591.51936
$ws.Range("I3").Value = 7
$ws.Range("J3").Value = "linearization_heuristic"
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 2231941.8

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1549.799349952191
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 43844049464.32219
$ws.Range("E4").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("F4").Value = 2000
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 40250806151.80891
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = "linearization_heuristic"
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 2231941.8

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 1577.369719139921
$ws.Range("C5").Value = 0.5
$ws.Range("D5").Value = 43845881124.41301
$ws.Range("E5").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("F5").Value = 2000
$ws.Range("G5").Value = 14
$ws.Range("H5").Value = 40185292430.46282
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = "linearization_heuristic"
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2231941.8

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1420.550188517038
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 43986747680.24268
$ws.Range("E6").Value = "linearization_heuristic"
$ws.Range("F6").Value = 2000
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 40703357020.42834
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = "linearization_heuristic"
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 2231941.8

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 1386.685491618477
$ws.Range("C7").Value = 0.5
$ws.Range("D7").Value = 43973311978.57693
$ws.Range("E7").Value = "linearization_heuristic"
$ws.Range("F7").Value = 2000
$ws.Range("G7").Value = 14
$ws.Range("H7").Value = 40771049891.36427
$ws.Range("I7").Value = 7
$ws.Range("J7").Value = "linearization_heuristic"
$ws.Range("K7").Value = 30000
$ws.Range("L7").Value = 2231941.8

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1420.530868342983
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 43986747680.24716
$ws.Range("E8").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("F8").Value = 2000
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 40703404012.3964
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = "linearization_heuristic"
$ws.Range("K8").Value = 30000
$ws.Range("L8").Value = 2231941.8

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 1386.669162326891
$ws.Range("C9").Value = 0.5
$ws.Range("D9").Value = 43973311978.57693
$ws.Range("E9").Value = "linearization_heuristic_Prop_Bouncing"
$ws.Range("F9").Value = 2000
$ws.Range("G9").Value = 14
$ws.Range("H9").Value = 40771089914.89458
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = "linearization_heuristic"
$ws.Range("K9").Value = 30000
$ws.Range("L9").Value = 2231941.8

